$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing row (2061) down through the new rows (2062:2083)
$ws.Range("A2061:K2061").Copy() | Out-Null
$ws.Range("A2062:K2083").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 2062
$ws.Range("A2062").Value = 45901
$ws.Range("B2062").Value = 1908.12
$ws.Range("C2062").Value = 1924.78
$ws.Range("D2062").Value = 1903.72
$ws.Range("E2062").Value = 1916.16
$ws.Range("F2062").Formula = "=E2062/1000"
$ws.Range("G2062").Value = 340257528
$ws.Range("H2062").Value = 683140000000
$ws.Range("I2062").Formula = "=I2061+1"
$ws.Range("J2062").Formula = "=SUM(`$F`$3:F2062)/I2062"
$ws.Range("K2062").Formula = "=SUM(F1813:F2062)/250"

# Row 2063
$ws.Range("A2063").Value = 45902
$ws.Range("B2063").Value = 1913.85
$ws.Range("C2063").Value = 1913.85
$ws.Range("D2063").Value = 1847.93
$ws.Range("E2063").Value = 1864.49
$ws.Range("F2063").Formula = "=E2063/1000"
$ws.Range("G2063").Value = 374441649
$ws.Range("H2063").Value = 702928000000
$ws.Range("I2063").Formula = "=I2062+1"
$ws.Range("J2063").Formula = "=SUM(`$F`$3:F2063)/I2063"
$ws.Range("K2063").Formula = "=SUM(F1814:F2063)/250"

# Row 2064
$ws.Range("A2064").Value = 45903
$ws.Range("B2064").Value = 1868.23
$ws.Range("C2064").Value = 1874.61
$ws.Range("D2064").Value = 1817.22
$ws.Range("E2064").Value = 1822.59
$ws.Range("F2064").Formula = "=E2064/1000"
$ws.Range("G2064").Value = 305322894
$ws.Range("H2064").Value = 547639000000
$ws.Range("I2064").Formula = "=I2063+1"
$ws.Range("J2064").Formula = "=SUM(`$F`$3:F2064)/I2064"
$ws.Range("K2064").Formula = "=SUM(F1815:F2064)/250"

# Row 2065
$ws.Range("A2065").Value = 45904
$ws.Range("B2065").Value = 1824.58
$ws.Range("C2065").Value = 1836.38
$ws.Range("D2065").Value = 1770.01
$ws.Range("E2065").Value = 1797.09
$ws.Range("F2065").Formula = "=E2065/1000"
$ws.Range("G2065").Value = 306215413
$ws.Range("H2065").Value = 569476000000
$ws.Range("I2065").Formula = "=I2064+1"
$ws.Range("J2065").Formula = "=SUM(`$F`$3:F2065)/I2065"
$ws.Range("K2065").Formula = "=SUM(F1816:F2065)/250"

# Row 2066
$ws.Range("A2066").Value = 45905
$ws.Range("B2066").Value = 1802.5
$ws.Range("C2066").Value = 1838.81
$ws.Range("D2066").Value = 1783.9
$ws.Range("E2066").Value = 1838.81
$ws.Range("F2066").Formula = "=E2066/1000"
$ws.Range("G2066").Value = 266539082
$ws.Range("H2066").Value = 478488000000
$ws.Range("I2066").Formula = "=I2065+1"
$ws.Range("J2066").Formula = "=SUM(`$F`$3:F2066)/I2066"
$ws.Range("K2066").Formula = "=SUM(F1817:F2066)/250"

# Row 2067
$ws.Range("A2067").Value = 45908
$ws.Range("B2067").Value = 1837.93
$ws.Range("C2067").Value = 1862.56
$ws.Range("D2067").Value = 1836.3
$ws.Range("E2067").Value = 1860.99
$ws.Range("F2067").Formula = "=E2067/1000"
$ws.Range("G2067").Value = 271428141
$ws.Range("H2067").Value = 491426000000
$ws.Range("I2067").Formula = "=I2066+1"
$ws.Range("J2067").Formula = "=SUM(`$F`$3:F2067)/I2067"
$ws.Range("K2067").Formula = "=SUM(F1818:F2067)/250"

# Row 2068
$ws.Range("A2068").Value = 45909
$ws.Range("B2068").Value = 1858.18
$ws.Range("C2068").Value = 1858.18
$ws.Range("D2068").Value = 1823.51
$ws.Range("E2068").Value = 1831.47
$ws.Range("F2068").Formula = "=E2068/1000"
$ws.Range("G2068").Value = 253780153
$ws.Range("H2068").Value = 425497000000
$ws.Range("I2068").Formula = "=I2067+1"
$ws.Range("J2068").Formula = "=SUM(`$F`$3:F2068)/I2068"
$ws.Range("K2068").Formula = "=SUM(F1819:F2068)/250"

# Row 2069
$ws.Range("A2069").Value = 45910
$ws.Range("B2069").Value = 1832.71
$ws.Range("C2069").Value = 1853.35
$ws.Range("D2069").Value = 1832.13
$ws.Range("E2069").Value = 1844.1
$ws.Range("F2069").Formula = "=E2069/1000"
$ws.Range("G2069").Value = 254058039
$ws.Range("H2069").Value = 434218000000
$ws.Range("I2069").Formula = "=I2068+1"
$ws.Range("J2069").Formula = "=SUM(`$F`$3:F2069)/I2069"
$ws.Range("K2069").Formula = "=SUM(F1820:F2069)/250"

# Row 2070
$ws.Range("A2070").Value = 45911
$ws.Range("B2070").Value = 1843.53
$ws.Range("C2070").Value = 1881.38
$ws.Range("D2070").Value = 1826.26
$ws.Range("E2070").Value = 1881.38
$ws.Range("F2070").Formula = "=E2070/1000"
$ws.Range("G2070").Value = 312111432
$ws.Range("H2070").Value = 582720000000
$ws.Range("I2070").Formula = "=I2069+1"
$ws.Range("J2070").Formula = "=SUM(`$F`$3:F2070)/I2070"
$ws.Range("K2070").Formula = "=SUM(F1821:F2070)/250"

# Row 2071
$ws.Range("A2071").Value = 45912
$ws.Range("B2071").Value = 1884.21
$ws.Range("C2071").Value = 1897.91
$ws.Range("D2071").Value = 1874.9
$ws.Range("E2071").Value = 1881.63
$ws.Range("F2071").Formula = "=E2071/1000"
$ws.Range("G2071").Value = 313935966
$ws.Range("H2071").Value = 602304000000
$ws.Range("I2071").Formula = "=I2070+1"
$ws.Range("J2071").Formula = "=SUM(`$F`$3:F2071)/I2071"
$ws.Range("K2071").Formula = "=SUM(F1822:F2071)/250"

# Row 2072
$ws.Range("A2072").Value = 45915
$ws.Range("B2072").Value = 1881.25
$ws.Range("C2072").Value = 1884.54
$ws.Range("D2072").Value = 1868.26
$ws.Range("E2072").Value = 1877.23
$ws.Range("F2072").Formula = "=E2072/1000"
$ws.Range("G2072").Value = 271831119
$ws.Range("H2072").Value = 522140000000
$ws.Range("I2072").Formula = "=I2071+1"
$ws.Range("J2072").Formula = "=SUM(`$F`$3:F2072)/I2072"
$ws.Range("K2072").Formula = "=SUM(F1823:F2072)/250"

# Row 2073
$ws.Range("A2073").Value = 45916
$ws.Range("B2073").Value = 1877.81
$ws.Range("C2073").Value = 1905.05
$ws.Range("D2073").Value = 1873.61
$ws.Range("E2073").Value = 1905.05
$ws.Range("F2073").Formula = "=E2073/1000"
$ws.Range("G2073").Value = 307828945
$ws.Range("H2073").Value = 561404000000
$ws.Range("I2073").Formula = "=I2072+1"
$ws.Range("J2073").Formula = "=SUM(`$F`$3:F2073)/I2073"
$ws.Range("K2073").Formula = "=SUM(F1824:F2073)/250"

# Row 2074
$ws.Range("A2074").Value = 45917
$ws.Range("B2074").Value = 1901.52
$ws.Range("C2074").Value = 1915.12
$ws.Range("D2074").Value = 1894.34
$ws.Range("E2074").Value = 1910.64
$ws.Range("F2074").Formula = "=E2074/1000"
$ws.Range("G2074").Value = 308418344
$ws.Range("H2074").Value = 553689000000
$ws.Range("I2074").Formula = "=I2073+1"
$ws.Range("J2074").Formula = "=SUM(`$F`$3:F2074)/I2074"
$ws.Range("K2074").Formula = "=SUM(F1825:F2074)/250"

# Row 2075
$ws.Range("A2075").Value = 45918
$ws.Range("B2075").Value = 1909.36
$ws.Range("C2075").Value = 1932.08
$ws.Range("D2075").Value = 1865.36
$ws.Range("E2075").Value = 1884.32
$ws.Range("F2075").Formula = "=E2075/1000"
$ws.Range("G2075").Value = 398466315
$ws.Range("H2075").Value = 758011000000
$ws.Range("I2075").Formula = "=I2074+1"
$ws.Range("J2075").Formula = "=SUM(`$F`$3:F2075)/I2075"
$ws.Range("K2075").Formula = "=SUM(F1826:F2075)/250"

# Row 2076
$ws.Range("A2076").Value = 45919
$ws.Range("B2076").Value = 1884.04
$ws.Range("C2076").Value = 1895.66
$ws.Range("D2076").Value = 1862.09
$ws.Range("E2076").Value = 1870.6
$ws.Range("F2076").Formula = "=E2076/1000"
$ws.Range("G2076").Value = 284001680
$ws.Range("H2076").Value = 525645000000
$ws.Range("I2076").Formula = "=I2075+1"
$ws.Range("J2076").Formula = "=SUM(`$F`$3:F2076)/I2076"
$ws.Range("K2076").Formula = "=SUM(F1827:F2076)/250"

# Row 2077
$ws.Range("A2077").Value = 45922
$ws.Range("B2077").Value = 1873.11
$ws.Range("C2077").Value = 1882.8
$ws.Range("D2077").Value = 1865.49
$ws.Range("E2077").Value = 1882.8
$ws.Range("F2077").Formula = "=E2077/1000"
$ws.Range("G2077").Value = 266808413
$ws.Range("H2077").Value = 517291000000
$ws.Range("I2077").Formula = "=I2076+1"
$ws.Range("J2077").Formula = "=SUM(`$F`$3:F2077)/I2077"
$ws.Range("K2077").Formula = "=SUM(F1828:F2077)/250"

# Row 2078
$ws.Range("A2078").Value = 45923
$ws.Range("B2078").Value = 1880.68
$ws.Range("C2078").Value = 1880.68
$ws.Range("D2078").Value = 1811.48
$ws.Range("E2078").Value = 1848.27
$ws.Range("F2078").Formula = "=E2078/1000"
$ws.Range("G2078").Value = 319070074
$ws.Range("H2078").Value = 581984000000
$ws.Range("I2078").Formula = "=I2077+1"
$ws.Range("J2078").Formula = "=SUM(`$F`$3:F2078)/I2078"
$ws.Range("K2078").Formula = "=SUM(F1829:F2078)/250"

# Row 2079
$ws.Range("A2079").Value = 45924
$ws.Range("B2079").Value = 1839.58
$ws.Range("C2079").Value = 1880.92
$ws.Range("D2079").Value = 1830.53
$ws.Range("E2079").Value = 1880.91
$ws.Range("F2079").Formula = "=E2079/1000"
$ws.Range("G2079").Value = 277896091
$ws.Range("H2079").Value = 546717000000
$ws.Range("I2079").Formula = "=I2078+1"
$ws.Range("J2079").Formula = "=SUM(`$F`$3:F2079)/I2079"
$ws.Range("K2079").Formula = "=SUM(F1830:F2079)/250"

# Row 2080
$ws.Range("A2080").Value = 45925
$ws.Range("B2080").Value = 1879.1
$ws.Range("C2080").Value = 1899.92
$ws.Range("D2080").Value = 1876.65
$ws.Range("E2080").Value = 1878.93
$ws.Range("F2080").Formula = "=E2080/1000"
$ws.Range("G2080").Value = 279857181
$ws.Range("H2080").Value = 562362000000
$ws.Range("I2080").Formula = "=I2079+1"
$ws.Range("J2080").Formula = "=SUM(`$F`$3:F2080)/I2080"
$ws.Range("K2080").Formula = "=SUM(F1831:F2080)/250"

# Row 2081
$ws.Range("A2081").Value = 45926
$ws.Range("B2081").Value = 1872.5
$ws.Range("C2081").Value = 1873.42
$ws.Range("D2081").Value = 1845.85
$ws.Range("E2081").Value = 1845.85
$ws.Range("F2081").Formula = "=E2081/1000"
$ws.Range("G2081").Value = 254215815
$ws.Range("H2081").Value = 491514000000
$ws.Range("I2081").Formula = "=I2080+1"
$ws.Range("J2081").Formula = "=SUM(`$F`$3:F2081)/I2081"
$ws.Range("K2081").Formula = "=SUM(F1832:F2081)/250"

# Row 2082
$ws.Range("A2082").Value = 45929
$ws.Range("B2082").Value = 1847.36
$ws.Range("C2082").Value = 1863.99
$ws.Range("D2082").Value = 1829.39
$ws.Range("E2082").Value = 1859.76
$ws.Range("F2082").Formula = "=E2082/1000"
$ws.Range("G2082").Value = 234605202
$ws.Range("H2082").Value = 457333000000
$ws.Range("I2082").Formula = "=I2081+1"
$ws.Range("J2082").Formula = "=SUM(`$F`$3:F2082)/I2082"
$ws.Range("K2082").Formula = "=SUM(F1833:F2082)/250"

# Row 2083
$ws.Range("A2083").Value = 45930
$ws.Range("B2083").Value = 1868.19
$ws.Range("C2083").Value = 1880.79
$ws.Range("D2083").Value = 1868.02
$ws.Range("E2083").Value = 1873.85
$ws.Range("F2083").Formula = "=E2083/1000"
$ws.Range("G2083").Value = 227571997
$ws.Range("H2083").Value = 468250000000
$ws.Range("I2083").Formula = "=I2082+1"
$ws.Range("J2083").Formula = "=SUM(`$F`$3:F2083)/I2083"
$ws.Range("K2083").Formula = "=SUM(F1834:F2083)/250"

# Move the selection to mirror the post-edit cursor position (A2084),
# matching the recorded sheetView/selection in the target workbook.
$ws.Range("A2084").Select() | Out-Null

# Force a full recalculation so every formula (esp. the rolling-window
# SUM ranges that now include the newly appended rows) carries a fresh,
# non-stale cached value before the workbook is saved.
$excel.CalculateFull()

Write-Output "done"
